$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B/C columns with new values (row.names text stays same)
$ws.Range("B2").Value = 0.464699131312926
$ws.Range("C2").Value = 0.490331077891466

$ws.Range("B3").Value = 0.515997903617863
$ws.Range("C3").Value = 0.524163591347172

$ws.Range("B4").Value = 0.65527111860267
$ws.Range("C4").Value = 0.802733037540415

$ws.Range("B5").Value = 1.88861055097863
$ws.Range("C5").Value = 2.24928758578876

$ws.Range("B6").Value = 0.502708803448382
$ws.Range("C6").Value = 0.523339725627166

# Resize the table to include the new column first, then set header text
# (setting the header text before Resize leaves the ListColumn named "Column4")
$table = $ws.ListObjects.Item("Table3")
$table.Resize($ws.Range("A1:D6"))

# New column D: geiger.averaged
$ws.Range("D1").Value = "geiger.averaged"
$ws.Range("D2").Value = 1.59445809846249
$ws.Range("D3").Value = 0.683199800574798
$ws.Range("D6").Value = 4.25705506189245
